$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 2
$ws.Range("E2").Value = 35

# Row 4
$ws.Range("E4").Value = 13

# Row 9
$ws.Range("E9").Value = 17

# Row 15
$ws.Range("E15").Value = 126
$ws.Range("F15").Value = 63
$ws.Range("H15").Value = 63

# Row 18
$ws.Range("E18").Value = 80

# Row 19
$ws.Range("E19").Value = 33

# Row 33
$ws.Range("E33").Value = 22

# Row 36
$ws.Range("E36").Value = 65
$ws.Range("F36").Value = 25
$ws.Range("H36").Value = 25

# Row 37
$ws.Range("E37").Value = 32

# Row 39
$ws.Range("E39").Value = 17
$ws.Range("F39").Value = 11
$ws.Range("H39").Value = 11

# Row 43
$ws.Range("E43").Value = 17
$ws.Range("F43").Value = 11
$ws.Range("H43").Value = 11

# Row 50
$ws.Range("E50").Value = 14

# Row 76
$ws.Range("E76").Value = 32

# Row 78
$ws.Range("E78").Value = 29
$ws.Range("F78").Value = 11
$ws.Range("H78").Value = 11

# Row 79
$ws.Range("E79").Value = 21

# Row 82
$ws.Range("E82").Value = 7

# Row 83
$ws.Range("E83").Value = 8

# Row 87
$ws.Range("E87").Value = 10

# Row 88
$ws.Range("E88").Value = 13

# Row 89
$ws.Range("E89").Value = 21
$ws.Range("F89").Value = 10
$ws.Range("H89").Value = 10

$wb.Save()
